# Remove the now-relocated data block (L5:N7, with L5:N5 merged) from the
# "ZSL Prompt Tuning" sheet.
$wb = $excel.ActiveWorkbook
$zsl = $wb.Worksheets.Item("ZSL Prompt Tuning")
$zsl.Range("L5:N5").UnMerge()
$zsl.Range("L5:N7").Clear()
$zsl.Range("L5:N7").Select()

# Add the new "AMI Chunking" sheet and move it to the 2nd position
# (right after "ZSL Prompt Tuning").
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "AMI Chunking"

# Re-fetch "ZSL Prompt Tuning" since the $zsl reference can rebind after the
# sheet collection is mutated by Add().
$zslAfterAdd = $wb.Worksheets.Item("ZSL Prompt Tuning")
$newSheet.Move($null, $zslAfterAdd)

# Re-fetch the worksheet object by name since the reference used to add/move
# it can go stale once it has been repositioned.
$chunk = $wb.Worksheets.Item("AMI Chunking")

# Header row: Model / Quantization / Scenario + merged "AMI" results header
$chunk.Range("A1").Value = "Model"
$chunk.Range("B1").Value = "Quantization"
$chunk.Range("C1").Value = "Scenario"
$chunk.Range("A1:C1").Font.Bold = $true
$chunk.Range("D1").Value = "AMI"
$chunk.Range("D1:F1").Merge()

# Sub-header row with metric names
$chunk.Range("D2").Value = "ROUGE-2"
$chunk.Range("E2").Value = "ROUGE-L"
$chunk.Range("F2").Value = "BERT-Score"

# Data that used to live in the "ZSL Prompt Tuning" sheet (L5:N7) for the
# xgen-7b-inst / 2SL scenario.
$chunk.Range("A3").Value = "xgen-7b-inst"
$chunk.Range("B3").Value = "4bit full"
$chunk.Range("C3").Value = "2SL"
$chunk.Range("D3").Value = 0.065877807339252234
$chunk.Range("E3").Value = 0.12999817316708137
$chunk.Range("F3").Value = 0.17726770726342958

# New placeholder rows for additional models/scenarios.
$chunk.Range("A4").Value = "xgen-7b-inst"
$chunk.Range("B4").Value = "4bit full"
$chunk.Range("C4").Value = "0SL"

$chunk.Range("A5").Value = "falcon-7b-inst"
$chunk.Range("B5").Value = "4bit full"
$chunk.Range("C5").Value = "0SL"

# Column widths matching the target layout.
$chunk.Columns.Item(1).ColumnWidth = 15.7109375
$chunk.Columns.Item(2).ColumnWidth = 14.28515625

# Make this the active sheet / selected tab (was previously "AMI 2SL C1").
$chunk.Activate()
$chunk.Range("D4").Select()
